# informe de 02-11-2020 al 07-11-2020
# Adds the week of Mon 2020-11-02 .. Sun 2020-11-08 to the "INFORME OCTUBRE"
# weekly report sheet: a date header row (24), a description row (25) and a
# trailing blank row (26); also grows the Tabla43 listobject to cover them.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("INFORME OCTUBRE")

# ---------------------------------------------------------------------------
# Row 24 - date header (Mon 2020-11-02 .. Sun 2020-11-08), same look as the
# other date-header rows (e.g. row 20/22).
# ---------------------------------------------------------------------------
$ws.Range("A20:H20").Copy()
$ws.Range("A24:H24").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A24").Value = 44137
$ws.Range("B24").Value = 44138
$ws.Range("C24").Value = 44139
$ws.Range("D24").Value = 44140
$ws.Range("E24").Value = 44141
$ws.Range("F24").Value = 44142
$ws.Range("G24").Value = 44143

# ---------------------------------------------------------------------------
# Row 25 - weekly description (developer notes for that week).
# Columns A, E, F use the wrap-text description style (same as C15 etc.);
# columns B, C, D, G use the plain style (same as F23); H uses the
# developer-name style, vertically centered.
# ---------------------------------------------------------------------------
$ws.Range("C15").Copy()
$ws.Range("A25").PasteSpecial(-4122)
$ws.Range("E25").PasteSpecial(-4122)
$ws.Range("F25").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("F23").Copy()
$ws.Range("B25").PasteSpecial(-4122)
$ws.Range("C25").PasteSpecial(-4122)
$ws.Range("D25").PasteSpecial(-4122)
$ws.Range("G25").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A25").Value = "Formulario Calculo( metodos para mostrar empleados ), modificacion tablas en SQLSERVER y codigo en C#(tipo de documento, contrato, regimen salud) "
$ws.Range("C25").Value = "-"
$ws.Range("D25").Value = "-"
$ws.Range("E25").Value = "Modificacion formulario calculo( diseño de agregar y quitar conceptos de la tabla calculo), Cambiar diseño  de calculo de planilla y inicio de sesion( según acordado)"
$ws.Range("F25").Value = "diseño inicio de sesion y revision de codigo en SQLSERVER (se borro los procedimientos que ya no se usaran.(según acordado)"

$ws.Range("H25").Value = "CARLOS MEZA"
$ws.Range("H25").VerticalAlignment = -4108   # xlCenter

$ws.Rows.Item(25).RowHeight = 135

# ---------------------------------------------------------------------------
# Row 26 - trailing blank spacer row (same plain style as F23 / the rest of
# row 25's blank cells). Only columns A:G are touched, H is left untouched.
# ---------------------------------------------------------------------------
$ws.Range("F23").Copy()
$ws.Range("A26").PasteSpecial(-4122)
$ws.Range("B26").PasteSpecial(-4122)
$ws.Range("C26").PasteSpecial(-4122)
$ws.Range("D26").PasteSpecial(-4122)
$ws.Range("E26").PasteSpecial(-4122)
$ws.Range("F26").PasteSpecial(-4122)
$ws.Range("G26").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Grow the "Tabla43" table (and its autofilter) so it keeps covering the
# data, now through row 26.
# ---------------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A9:H26"))

# ---------------------------------------------------------------------------
# Selection / scroll position, matching where the author left the sheet.
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("D25").Select()
